# Update the "No. of pmid" and "No. of doi" sheets with the new MeSH table values.
$wb = $excel.ActiveWorkbook

# Column headers B..Z correspond to years 2000..2024; rows 2-4 hold the three concepts.
$columns = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")

# New values per sheet, per row (2 = Artificial intelligence, 3 = Machine learning, 4 = Natural language processing)
$newValues = @{
    "No. of pmid" = @{
        2 = @(7959,8482,9041,9883,10819,12091,13464,14584,14865,16335,16500,17759,18646,19905,20291,20632,20406,20888,22356,24926,27756,32815,34011,25922,197)
        3 = @(1562,1850,1892,2193,2517,2732,3158,3573,3880,4174,4492,5009,5167,5949,6257,6563,6849,7204,8190,9575,11099,13970,15518,11591,106)
        4 = @(182,177,224,238,261,291,315,303,301,286,343,371,412,427,378,424,409,413,481,534,534,628,697,531,1)
    }
    "No. of doi" = @{
        2 = @(6451,6925,7473,8341,9323,10468,11931,13110,13592,15064,15408,16589,17482,18772,19161,19527,19571,20288,21906,24605,27477,32509,33779,25782,197)
        3 = @(1287,1534,1581,1878,2205,2415,2841,3232,3547,3883,4198,4730,4852,5607,5964,6232,6576,7023,8031,9439,11003,13863,15446,11545,106)
        4 = @(152,148,173,181,204,225,251,237,235,257,298,332,365,389,338,381,375,373,445,521,522,616,685,532,1)
    }
}

foreach ($sheetName in $newValues.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsData = $newValues[$sheetName]
    foreach ($rowNum in $rowsData.Keys) {
        $values = $rowsData[$rowNum]
        for ($i = 0; $i -lt $columns.Length; $i++) {
            $cellRef = "$($columns[$i])$rowNum"
            $ws.Range($cellRef).Value = $values[$i]
        }
    }
}
